# Optimize Streamlit startup and lazy-load models
# Updates the Random Forest RMSE metric on the Summary sheet and refreshes
# the forecasted-volume figures on the Gradient Boosting and Random Forest
# sheets to reflect the new lazy-loaded model outputs.

$wb = $excel.ActiveWorkbook

# --- Summary sheet: Random Forest row (row 5) ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("C5").Value = 165164.3121258343
$summary.Range("D5").Value = "Random Forest ranked 4th based on RMSE of 165164.31. It performs well for short-term, interpretable."

# --- Gradient Boosting sheet: forecasted volume values ---
# 11 new values cycle across B2:B12 and repeat across B14:B24 (B13/B25 keep
# their original, unchanged value).
$gbNew = @(20009580.75716664, 20505414.43236707, 20722176.28794238, 20239350.86255708, 20593523.26186129, 20874182.11893817, 21212990.68746854, 22048252.03013377, 22469992.57835477, 23380378.42595203, 24014981.72558681)

$gb = $wb.Worksheets.Item("Gradient Boosting")
for ($i = 0; $i -lt $gbNew.Length; $i++) {
    $gb.Cells.Item(2 + $i, 2).Value = $gbNew[$i]
    $gb.Cells.Item(14 + $i, 2).Value = $gbNew[$i]
}

# --- Random Forest sheet: forecasted volume values ---
# 12 new values cycle across B2:B13 and repeat across B14:B25.
$rfNew = @(20596000, 20664000, 20803000, 20719000, 20772000, 21148000, 21650000, 21959000, 22228000, 22889000, 23496000, 24100000)

$rf = $wb.Worksheets.Item("Random Forest")
for ($i = 0; $i -lt $rfNew.Length; $i++) {
    $rf.Cells.Item(2 + $i, 2).Value = $rfNew[$i]
    $rf.Cells.Item(14 + $i, 2).Value = $rfNew[$i]
}
